# Update CDA Logical model for ST.r2b
# - Rename "Include ValueSets" -> "Include ValueSet #0"
# - Rename "Include from EntityCode" -> "Include #1"
# - Bump Version / Date on the Metadata sheet
# - Insert a new "Jurisdiction" property row (empty value) right after
#   "Contact" and before "Description" on the Metadata sheet

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsValueSets = $wb.Worksheets.Item("Include ValueSets")
$wsEntityCode = $wb.Worksheets.Item("Include from EntityCode")

# --- Sheet renames -------------------------------------------------------
$wsValueSets.Name = "Include ValueSet #0"
$wsEntityCode.Name = "Include #1"

# --- Metadata: bump Version + Date ---------------------------------------
$wsMetadata.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$wsMetadata.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- Metadata: insert "Jurisdiction" row (row 11, before "Description") --
$wsMetadata.Rows.Item(11).Insert()

# Copy formatting from the row that is now directly below (was row 11,
# "Description", now pushed to row 12) so the new row matches the existing
# property/value styling instead of the default "no border" style.
$wsMetadata.Range("A12:B12").Copy()
$wsMetadata.Range("A11:B11").PasteSpecial(-4122)

$wsMetadata.Range("A11").Value = "Jurisdiction"
$wsMetadata.Range("B11").Value = ""
